# Edit script for spain_laliga2_2023-2024 workbook
# 1) Swap the betting-row data (columns F:V) for several out-of-order row
#    pairs so that match rows line up chronologically again.
# 2) Append 9 new match rows (168-176) that were added at the end of the
#    sheet, copying the existing formatting used by the previous last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row pair swaps: for each pair, exchange the contents of columns
#    F..V (home team through match url). Columns A..E (Indice, pais,
#    torneio, temporada, data_partida) stay untouched.
# ---------------------------------------------------------------------
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$pairs = @(
    @(51,52),
    @(72,73),
    @(82,83),
    @(85,87),
    @(118,119),
    @(120,121),
    @(124,125),
    @(136,137)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

# ---------------------------------------------------------------------
# 2) Append new rows 168-176.
# ---------------------------------------------------------------------
$newRows = @{
    168 = @('167','spain','laliga2','2023-2024','45248.67708333334','Alcorcon','0','Gijon','0','2.89','12/11/2023 14:12','3.34','18/11/2023 16:14','3.04','12/11/2023 14:12','3.07','18/11/2023 16:13','2.75','12/11/2023 14:12','2.48','18/11/2023 16:14','https://www.betexplorer.com/football/spain/laliga2/alcorcon-gijon/tKucgIpC/')
    169 = @('168','spain','laliga2','2023-2024','45248.77083333334','Amorebieta','2','Tenerife','0','3.76','12/11/2023 18:42','2.92','18/11/2023 18:29','3.16','12/11/2023 18:42','2.88','18/11/2023 18:26','2.23','12/11/2023 18:42','2.96','18/11/2023 18:29','https://www.betexplorer.com/football/spain/laliga2/amorebieta-tenerife/OrV1hx0I/')
    170 = @('169','spain','laliga2','2023-2024','45248.77083333334','Zaragoza','0','Huesca','2','1.91','12/11/2023 18:42','2.14','18/11/2023 18:29','3.32','12/11/2023 18:42','3.03','18/11/2023 18:29','4.81','12/11/2023 18:42','4.34','18/11/2023 18:29','https://www.betexplorer.com/football/spain/laliga2/zaragoza-huesca/dfxYfVLG/')
    171 = @('170','spain','laliga2','2023-2024','45248.875','Espanyol','2','Elche','0','1.95','12/11/2023 18:42','1.86','18/11/2023 20:55','3.58','12/11/2023 18:42','3.58','18/11/2023 20:55','4.15','12/11/2023 18:42','4.71','18/11/2023 20:55','https://www.betexplorer.com/football/spain/laliga2/espanyol-elche/I3wUeB6A/')
    172 = @('171','spain','laliga2','2023-2024','45249.58333333334','Ferrol','1','Burgos CF','1','1.99','12/11/2023 14:12','2.1','19/11/2023 13:53','3.31','12/11/2023 14:12','3.11','19/11/2023 13:53','4.41','12/11/2023 14:12','4.33','19/11/2023 13:53','https://www.betexplorer.com/football/spain/laliga2/ferrol-burgos-cf/AwzMcXyb/')
    173 = @('172','spain','laliga2','2023-2024','45249.67708333334','Eldense','2','Mirandes','2','1.99','12/11/2023 16:43','2.23','19/11/2023 16:13','3.46','12/11/2023 16:43','3.26','19/11/2023 16:11','4.02','12/11/2023 16:43','3.66','19/11/2023 16:13','https://www.betexplorer.com/football/spain/laliga2/eldense-mirandes/08hVzAUp/')
    174 = @('173','spain','laliga2','2023-2024','45249.67708333334','R. Oviedo','2','Eibar','1','2.94','13/11/2023 22:12','2.76','19/11/2023 15:49','2.94','13/11/2023 22:12','2.87','19/11/2023 16:13','2.79','13/11/2023 22:12','3.1','19/11/2023 15:49','https://www.betexplorer.com/football/spain/laliga2/r-oviedo-eibar/OtvQdii4/')
    175 = @('174','spain','laliga2','2023-2024','45249.77083333334','FC Cartagena SAD','1','Albacete','1','2.87','13/11/2023 22:12','3.69','19/11/2023 18:27','3.15','13/11/2023 22:12','3.23','19/11/2023 18:24','2.68','13/11/2023 22:12','2.24','19/11/2023 18:27','https://www.betexplorer.com/football/spain/laliga2/fc-cartagena-sad-albacete/SEpHbDMi/')
    176 = @('175','spain','laliga2','2023-2024','45249.77083333334','Villarreal B','0','Andorra','0','2.54','12/11/2023 18:42','2.61','19/11/2023 18:29','3.3','12/11/2023 18:42','3.46','19/11/2023 18:04','2.99','12/11/2023 18:42','2.81','19/11/2023 18:29','https://www.betexplorer.com/football/spain/laliga2/villarreal-fc-andorra/xlG9MPu9/')
}

# Columns A..V, in order, matching the 22 values stored for each row above.
$allCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
# Numeric columns must be written as real numbers (not text).
$numericCols = @("A","E","G","I","J","L","N","P","R","T")

$lastDataRow = 167

$rowNums = $newRows.Keys | Sort-Object
foreach ($r in $rowNums) {
    $values = $newRows[$r]

    # Copy formatting for the "Indice" (A) and "data_partida" (E) columns
    # from the previous last row, since those are the only two columns
    # that carry an explicit style in this sheet.
    $ws.Range("A$lastDataRow").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("E$lastDataRow").Copy() | Out-Null
    $ws.Range("E$r").PasteSpecial(-4122) | Out-Null

    for ($i = 0; $i -lt $allCols.Length; $i++) {
        $col = $allCols[$i]
        $addr = "$col$r"
        $val = $values[$i]
        if ($numericCols -contains $col) {
            $ws.Range($addr).Value = [double]$val
        }
        else {
            $ws.Range($addr).Value = $val
        }
    }
}

$excel.CutCopyMode = 0
